# GoodInfo_v2 - 2022-01-19 未完成
# Roll the report forward by one trading day:
#   - shift the "AD" (latest) volume column into "AE" (previous day),
#     and populate "AD" with the new 01/19 volume figures
#   - relabel the date-stamped headers (AD1/AE1/AG1)
#   - refresh the per-stock metrics that changed with the new day's data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1 headers ---------------------------------------------------
$ws.Range("AD1").Value = "01/19成交量(張)"
$ws.Range("AE1").Value = "01/18成交量(張)"
$ws.Range("AG1").Value = "類型_2022-01-19"

# ---- Row 2 (2368 金像電) ----------------------------------------------
$ws.Range("L2").Value = 88.7
$ws.Range("M2").Value = 88.7
$ws.Range("N2").Value = 88.7

$ws.Range("Z2").Value = 19.3
$ws.Range("AA2").Value = 18.6
$ws.Range("AB2").Value = 26
$ws.Range("AC2").Value = 26.4

# AD/AE hold text like "49,471" - force text format so Excel doesn't
# reinterpret the comma-separated digits as a number.
$ws.Range("AE2").NumberFormat = "@"
$ws.Range("AE2").Value = "107,651"
$ws.Range("AD2").NumberFormat = "@"
$ws.Range("AD2").Value = "49,471"

$ws.Range("AO2").Value = $false
$ws.Range("AP2").Value = $false

# ---- Row 3 (3036 文曄) -------------------------------------------------
$ws.Range("Z3").Value = 0.53
$ws.Range("AA3").Value = -0.07000000000000001
$ws.Range("AB3").Value = 0.86
$ws.Range("AC3").Value = 0.49

$ws.Range("AE3").NumberFormat = "@"
$ws.Range("AE3").Value = "10,754"
$ws.Range("AD3").NumberFormat = "@"
$ws.Range("AD3").Value = "7,328"

$ws.Range("AM3").Value = 654

# ---- Row 4 (6235 華孚) -------------------------------------------------
$ws.Range("Z4").Value = 1.53
$ws.Range("AA4").Value = 0.41
$ws.Range("AB4").Value = -1.01
$ws.Range("AC4").Value = -3.08

$ws.Range("AE4").NumberFormat = "@"
$ws.Range("AE4").Value = "25,882"
$ws.Range("AD4").NumberFormat = "@"
$ws.Range("AD4").Value = "25,397"
